# Weekly data update: insert 4 new "Pera" price rows (Abate Fettel / Packham's
# Triumph, Primera/Segunda) for the latest reporting date, pushing the
# existing historical rows down by 4 (dimension grows from A1:T430 to
# A1:T434).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing rows 407:430 down to 411:434, inserting 4 blank rows at 407.
$ws.Range("A407:T410").Insert()

# Shared values for this market/product block (constant across all rows).
$mercadoId = 11
$mercado   = "Vega Monumental Concepción"
$region    = "Bíobío"
$codreg    = 8
$tipo      = "Fruta"
$prodId    = 100104
$producto  = "Frutos de pepita"
$catId     = 100104005
$categoria = "Pera"
$unidad    = "$/caja 16 kilos empedrada"
$kgUnidad  = 16

$newRows = @(
    @{ Row=407; Fecha=44714; Variedad="Abate Fettel";        Calidad="Primera"; Volumen=50; PMin=8000; PMax=8000; PProm=8000; Origen="Región de O'Higgins"; PKg=500 },
    @{ Row=408; Fecha=44714; Variedad="Abate Fettel";        Calidad="Segunda"; Volumen=50; PMin=7000; PMax=7000; PProm=7000; Origen="Región de O'Higgins"; PKg=438 },
    @{ Row=409; Fecha=44714; Variedad="Packham's Triumph";   Calidad="Primera"; Volumen=50; PMin=8000; PMax=8000; PProm=8000; Origen="Región de O'Higgins"; PKg=500 },
    @{ Row=410; Fecha=44714; Variedad="Packham's Triumph";   Calidad="Segunda"; Volumen=50; PMin=7000; PMax=7000; PProm=7000; Origen="Región de O'Higgins"; PKg=438 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value  = $mercadoId
    $ws.Cells.Item($row, 2).Value  = $mercado
    $ws.Cells.Item($row, 3).Value  = $region
    $ws.Cells.Item($row, 4).Value  = $r.Fecha
    $ws.Cells.Item($row, 5).Value  = $codreg
    $ws.Cells.Item($row, 6).Value  = $tipo
    $ws.Cells.Item($row, 7).Value  = $prodId
    $ws.Cells.Item($row, 8).Value  = $producto
    $ws.Cells.Item($row, 9).Value  = $catId
    $ws.Cells.Item($row, 10).Value = $categoria
    $ws.Cells.Item($row, 11).Value = $r.Variedad
    $ws.Cells.Item($row, 12).Value = $r.Calidad
    $ws.Cells.Item($row, 13).Value = $r.Volumen
    $ws.Cells.Item($row, 14).Value = $r.PMin
    $ws.Cells.Item($row, 15).Value = $r.PMax
    $ws.Cells.Item($row, 16).Value = $r.PProm
    $ws.Cells.Item($row, 17).Value = $unidad
    $ws.Cells.Item($row, 18).Value = $r.Origen
    $ws.Cells.Item($row, 19).Value = $r.PKg
    $ws.Cells.Item($row, 20).Value = $kgUnidad
}
